# Latest generated outputs 2025-10-28
# For rows 2-20 (the "Application" top-level group on the Specification sheet),
# insert a new leading field-column value "Application" into column C,
# shifting the existing C/D/E values along into D/E/F respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

for ($r = 2; $r -le 20; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    # Shift existing values right by one column: C->D, D->E, E->F
    $ws.Cells.Item($r, 6).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal

    # Insert the new "Application" label into column C
    $ws.Cells.Item($r, 3).Value = "Application"
}
